$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (B1:D1)
$ws.Range("B1").Value = "file"
$ws.Range("C1").Value = "uploaded_on"
$ws.Range("D1").Value = "description"

# Row 2: A2 is numeric 0, B2 must stay the literal text "0"
$ws.Range("A2").Value = 0

$helper = $ws.Range("Z99")
$helper.Value = "'0"
$helper.Copy()
$ws.Range("B2").PasteSpecial(-4163)
$helper.Clear()

# Build the header style (bold, centered/top aligned, thin box border) on B1 once,
# then stamp it onto the remaining cells via a format-only paste so no throw-away
# intermediate cell styles get left behind in the style table.
$r = $ws.Range("B1")
$r.Font.Bold = $true
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4160
$r.Borders.LineStyle = 1

$r.Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("A2").PasteSpecial(-4122)
